$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.979.59'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '3.410.33'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '409.37'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.53'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.633'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +6.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.734'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +5.59%  '
$ws.Range('E10').Value = '  +1.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.82'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000221'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +42.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.35'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +10.67%  '
$ws.Range('E14').Value = '  -0.35%  '
$ws.Range('D15').Value = '3.955.96'
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.23'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +6.92%  '
$ws.Range('D17').Value = '3.413.23'
$ws.Range('E17').Value = '  -0.94%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.52'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +8.77%  '
$ws.Range('E19').Value = '  +7.10%  '
$ws.Range('D20').Value = '62.018.13'
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '454.19'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +44.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '92.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +9.03%  '
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.29'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +3.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.34'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +14.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.18'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +11.58%  '
$ws.Range('E28').Value = '  +0.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.70'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.76'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.01'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.171'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.31%  '
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '42.91'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -3.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +4.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '53.73'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +3.67%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('E40').Value = '  +7.66%  '
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('E42').Value = '  -1.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '142.95'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.25'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +8.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.58'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +16.11%  '
$ws.Range('E46').Value = '  +1.15%  '
$ws.Range('E47').Value = '  -1.18%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.45'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +5.60%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.147'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +21.86%  '
$ws.Range('E50').Value = '  +9.51%  '
$ws.Range('D51').Value = '3.758.31'
$ws.Range('E51').Value = '  -0.81%  '
